$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2280.7334
$ws.Range("I94").Value = 2280.7334
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2280.7334
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1829.7334
$ws.Range("N94").ClearContents()

$ws.Range("H132").Value = 2599617
$ws.Range("I132").Value = 2749325.5
$ws.Range("J132").Value = 4668.6665
$ws.Range("K132").Value = 8247976.5
$ws.Range("L132").Value = 14005.9995
$ws.Range("M132").Value = -8245446.5
$ws.Range("N132").Value = -19065.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 867.2917
$ws.Range("I2").Value = 600.1429000000001
$ws.Range("J2").Value = 1241.3
$ws.Range("K2").Value = 600.1429000000001
$ws.Range("L2").Value = 1241.3
$ws.Range("M2").Value = -487.1429000000001
$ws.Range("N2").Value = -1467.3

$ws.Range("H32").Value = 22431.328
$ws.Range("I32").Value = 26573.355
$ws.Range("J32").Value = 3792.2
$ws.Range("K32").Value = 26573.355
$ws.Range("L32").Value = 3792.2
$ws.Range("M32").Value = -26286.355
$ws.Range("N32").Value = -4366.2

$ws.Range("H116").Value = 867.2917
$ws.Range("I116").Value = 600.1429000000001
$ws.Range("J116").Value = 1241.3
$ws.Range("K116").Value = 600.1429000000001
$ws.Range("L116").Value = 1241.3
$ws.Range("M116").Value = 1693.8571
$ws.Range("N116").Value = -5829.3

$ws.Range("H132").Value = 4010.4807
$ws.Range("I132").Value = 3728.1316
$ws.Range("K132").Value = 11184.3948
$ws.Range("M132").Value = -8654.3948

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 867.2917
$ws.Range("I3").Value = 600.1429000000001
$ws.Range("J3").Value = 1241.3
$ws.Range("K3").Value = 600.1429000000001
$ws.Range("L3").Value = 1241.3
$ws.Range("M3").Value = -486.1429000000001
$ws.Range("N3").Value = -1469.3

$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("M89").ClearContents()

$ws.Range("H99").Value = 624.2105
$ws.Range("I99").Value = 521.875
$ws.Range("J99").Value = 1170
$ws.Range("K99").Value = 521.875
$ws.Range("L99").Value = 1170
$ws.Range("M99").Value = 976.125
$ws.Range("N99").Value = -4166

$ws.Range("H105").Value = 2062.8928
$ws.Range("I105").Value = 1987.8572
$ws.Range("J105").Value = 2137.9285
$ws.Range("K105").Value = 1987.8572
$ws.Range("L105").Value = 2137.9285
$ws.Range("M105").Value = -240.8571999999999
$ws.Range("N105").Value = -5631.9285

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 5450
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5450
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5450
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5798

$ws.Range("H41").Value = 6250
$ws.Range("I41").Value = 2000
$ws.Range("J41").Value = 10500
$ws.Range("K41").Value = 2000
$ws.Range("L41").Value = 10500
$ws.Range("M41").Value = -1572
$ws.Range("N41").Value = -11356

$ws.Range("H50").Value = 15697.333
$ws.Range("J50").Value = 15697.333
$ws.Range("L50").Value = 15697.333
$ws.Range("N50").Value = -16947.333

$ws.Range("H51").Value = 10343.091
$ws.Range("J51").Value = 10343.091
$ws.Range("L51").Value = 10343.091
$ws.Range("N51").Value = -11815.091

$ws.Range("H59").Value = 30812.7
$ws.Range("I59").Value = 15000
$ws.Range("J59").Value = 32569.666
$ws.Range("K59").Value = 15000
$ws.Range("L59").Value = 32569.666
$ws.Range("M59").Value = -13855
$ws.Range("N59").Value = -34859.666

$ws.Range("H60").Value = 11743.909
$ws.Range("I60").Value = 3600
$ws.Range("J60").Value = 18530.5
$ws.Range("K60").Value = 3600
$ws.Range("L60").Value = 18530.5
$ws.Range("M60").Value = -3089
$ws.Range("N60").Value = -19552.5

$ws.Range("H61").Value = 10343.091
$ws.Range("J61").Value = 10343.091
$ws.Range("L61").Value = 10343.091
$ws.Range("N61").Value = -11039.091

$ws.Range("H94").Value = 2642.95
$ws.Range("I94").Value = 1772.1666
$ws.Range("J94").Value = 3016.1428
$ws.Range("K94").Value = 1772.1666
$ws.Range("L94").Value = 3016.1428
$ws.Range("M94").Value = -1321.1666
$ws.Range("N94").Value = -3918.1428

$ws.Range("H99").Value = 2551.3333
$ws.Range("I99").Value = 2113.7778
$ws.Range("J99").Value = 2988.889
$ws.Range("K99").Value = 2113.7778
$ws.Range("L99").Value = 2988.889
$ws.Range("M99").Value = -615.7777999999998
$ws.Range("N99").Value = -5984.889

$ws.Range("H126").Value = 2551.3333
$ws.Range("I126").Value = 2113.7778
$ws.Range("J126").Value = 2988.889
$ws.Range("K126").Value = 6341.3334
$ws.Range("L126").Value = 8966.667000000001
$ws.Range("M126").Value = -3871.3334
$ws.Range("N126").Value = -13906.667

$ws.Range("H134").Value = 905.2222
$ws.Range("I134").Value = 621
$ws.Range("K134").Value = 1863
$ws.Range("M134").Value = 672

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1285.41
$ws.Range("I68").Value = 1157.1837
$ws.Range("J68").Value = 1408.6078
$ws.Range("K68").Value = 3471.5511
$ws.Range("L68").Value = 4225.8234
$ws.Range("M68").Value = -2660.5511
$ws.Range("N68").Value = -5847.8234

$ws.Range("H71").Value = 1285.41
$ws.Range("I71").Value = 1157.1837
$ws.Range("J71").Value = 1408.6078
$ws.Range("K71").Value = 10414.6533
$ws.Range("L71").Value = 12677.4702
$ws.Range("M71").Value = -6358.6533
$ws.Range("N71").Value = -20789.4702

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 41000
$ws.Range("J20").Value = 41000
$ws.Range("L20").Value = 41000
$ws.Range("N20").Value = -41490

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("M24").ClearContents()

$ws.Range("H80").Value = 4326.2104
$ws.Range("I80").Value = 2542.2727
$ws.Range("J80").Value = 6779.125
$ws.Range("K80").Value = 2542.2727
$ws.Range("L80").Value = 6779.125
$ws.Range("M80").Value = -1544.2727
$ws.Range("N80").Value = -8775.125

$ws.Range("H83").Value = 4326.2104
$ws.Range("I83").Value = 2542.2727
$ws.Range("J83").Value = 6779.125
$ws.Range("K83").Value = 12711.3635
$ws.Range("L83").Value = 33895.625
$ws.Range("M83").Value = -7719.363499999999
$ws.Range("N83").Value = -43879.625

$ws.Range("H102").Value = 1086.1666
$ws.Range("I102").Value = 900.8
$ws.Range("K102").Value = 900.8
$ws.Range("M102").Value = 721.2

$ws.Range("H132").Value = 59386.855
$ws.Range("I132").Value = 92564.87
$ws.Range("J132").Value = 3239.4614
$ws.Range("K132").Value = 277694.61
$ws.Range("L132").Value = 9718.3842
$ws.Range("M132").Value = -275164.61
$ws.Range("N132").Value = -14778.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2024.8077
$ws.Range("I46").Value = 1709.1818
$ws.Range("J46").Value = 2256.2666
$ws.Range("K46").Value = 1709.1818
$ws.Range("L46").Value = 2256.2666
$ws.Range("M46").Value = -1521.1818
$ws.Range("N46").Value = -2632.2666

$ws.Range("H68").Value = 1026.8182
$ws.Range("I68").Value = 1026.8182
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1026.8182
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -277.8181999999999
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 1026.8182
$ws.Range("I71").Value = 1026.8182
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 5134.090999999999
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -1390.090999999999
$ws.Range("N71").ClearContents()

$ws.Range("H82").Value = 1612.5
$ws.Range("J82").Value = 816.6667
$ws.Range("L82").Value = 816.6667
$ws.Range("N82").Value = -1538.6667

$ws.Range("H85").Value = 1612.5
$ws.Range("J85").Value = 816.6667
$ws.Range("L85").Value = 816.6667
$ws.Range("N85").Value = -3312.6667

$ws.Range("H132").Value = 6331.614
$ws.Range("I132").Value = 10345.546
$ws.Range("J132").Value = 2317.682
$ws.Range("K132").Value = 31036.638
$ws.Range("L132").Value = 6953.045999999999
$ws.Range("M132").Value = -28506.638
$ws.Range("N132").Value = -12013.046

$ws.Range("H136").Value = 7403.95
$ws.Range("I136").Value = 10383.25
$ws.Range("J136").Value = 2935
$ws.Range("K136").Value = 31149.75
$ws.Range("L136").Value = 8805
$ws.Range("M136").Value = -28599.75
$ws.Range("N136").Value = -13905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2050
$ws.Range("I62").Value = 600
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 600
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = 24
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 2050
$ws.Range("I65").Value = 600
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 3000
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = 120
$ws.Range("N65").Value = -23740

$ws.Range("H126").Value = 1149.7667
$ws.Range("J126").Value = 663.4211
$ws.Range("L126").Value = 1990.2633
$ws.Range("N126").Value = -6930.263300000001

$ws.Range("H132").Value = 1439.5714
$ws.Range("I132").Value = 1124.5714
$ws.Range("K132").Value = 3373.7142
$ws.Range("M132").Value = -843.7142000000003
